# Updated cryptos list — refresh Price/Volume(1h) columns (and, for rows
# 49-50, the Coin/Link identity too, since Aave and RenderToken swapped
# ranking positions) to match the latest scrape.
#
# Column D ("Price") values are stored as plain text in this workbook even
# when they look numeric (e.g. "0.9991", "22.91") — that preserves exact
# decimal formatting (trailing zeros, tiny magnitudes like
# "0.00000000135") the way the original scraper wrote them. Excel's COM
# layer auto-converts a numeric-looking string assigned via .Value into a
# real number (dropping/altering that formatting), so for those cells we
# briefly force the cell to Text number-format before the assignment, then
# clear the format again so the cell's style stays the same as every other
# data cell (only the stored text differs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = 'D2'; Value = '29.072.34'; Mode = 'AsIs' },
    @{ Cell = 'E2'; Value = '  -1.82%  '; Mode = 'AsIs' },
    @{ Cell = 'D3'; Value = '1.834.05'; Mode = 'AsIs' },
    @{ Cell = 'E3'; Value = '  -1.33%  '; Mode = 'AsIs' },
    @{ Cell = 'D4'; Value = '0.9991'; Mode = 'AsText' },
    @{ Cell = 'E4'; Value = '  +0.02%  '; Mode = 'AsIs' },
    @{ Cell = 'D5'; Value = '239.97'; Mode = 'AsText' },
    @{ Cell = 'E5'; Value = '  -1.97%  '; Mode = 'AsIs' },
    @{ Cell = 'D6'; Value = '0.6718'; Mode = 'AsText' },
    @{ Cell = 'E6'; Value = '  -3.23%  '; Mode = 'AsIs' },
    @{ Cell = 'E7'; Value = '  +0.06%  '; Mode = 'AsIs' },
    @{ Cell = 'D8'; Value = '0.07431'; Mode = 'AsText' },
    @{ Cell = 'E8'; Value = '  -3.45%  '; Mode = 'AsIs' },
    @{ Cell = 'D9'; Value = '0.2964'; Mode = 'AsText' },
    @{ Cell = 'E9'; Value = '  -3.11%  '; Mode = 'AsIs' },
    @{ Cell = 'D10'; Value = '22.91'; Mode = 'AsText' },
    @{ Cell = 'E10'; Value = '  -3.45%  '; Mode = 'AsIs' },
    @{ Cell = 'D11'; Value = '0.07656'; Mode = 'AsText' },
    @{ Cell = 'E11'; Value = '  -1.54%  '; Mode = 'AsIs' },
    @{ Cell = 'D12'; Value = '1.832.96'; Mode = 'AsIs' },
    @{ Cell = 'E12'; Value = '  -1.24%  '; Mode = 'AsIs' },
    @{ Cell = 'D13'; Value = '5.003'; Mode = 'AsText' },
    @{ Cell = 'E13'; Value = '  -2.78%  '; Mode = 'AsIs' },
    @{ Cell = 'D14'; Value = '0.6742'; Mode = 'AsText' },
    @{ Cell = 'E14'; Value = '  -2.66%  '; Mode = 'AsIs' },
    @{ Cell = 'D15'; Value = '86.37'; Mode = 'AsText' },
    @{ Cell = 'E15'; Value = '  -5.62%  '; Mode = 'AsIs' },
    @{ Cell = 'D16'; Value = '6.133'; Mode = 'AsText' },
    @{ Cell = 'E16'; Value = '  -6.64%  '; Mode = 'AsIs' },
    @{ Cell = 'D17'; Value = '29.064.45'; Mode = 'AsIs' },
    @{ Cell = 'E17'; Value = '  -1.77%  '; Mode = 'AsIs' },
    @{ Cell = 'D18'; Value = '0.000008236'; Mode = 'AsText' },
    @{ Cell = 'E18'; Value = '  -0.68%  '; Mode = 'AsIs' },
    @{ Cell = 'D19'; Value = '227.24'; Mode = 'AsText' },
    @{ Cell = 'D20'; Value = '12.46'; Mode = 'AsText' },
    @{ Cell = 'E20'; Value = '  -2.49%  '; Mode = 'AsIs' },
    @{ Cell = 'D21'; Value = '0.9994'; Mode = 'AsText' },
    @{ Cell = 'E21'; Value = '  -0.02%  '; Mode = 'AsIs' },
    @{ Cell = 'D22'; Value = '7.317'; Mode = 'AsText' },
    @{ Cell = 'E22'; Value = '  -3.79%  '; Mode = 'AsIs' },
    @{ Cell = 'D23'; Value = '0.9996'; Mode = 'AsText' },
    @{ Cell = 'E23'; Value = '  -0.03%  '; Mode = 'AsIs' },
    @{ Cell = 'D24'; Value = '160.38'; Mode = 'AsText' },
    @{ Cell = 'E24'; Value = '  +0.39%  '; Mode = 'AsIs' },
    @{ Cell = 'D25'; Value = '0.1430'; Mode = 'AsText' },
    @{ Cell = 'E25'; Value = '  -4.68%  '; Mode = 'AsIs' },
    @{ Cell = 'D26'; Value = '8.686'; Mode = 'AsText' },
    @{ Cell = 'E26'; Value = '  -2.72%  '; Mode = 'AsIs' },
    @{ Cell = 'D27'; Value = '17.97'; Mode = 'AsText' },
    @{ Cell = 'E27'; Value = '  -1.69%  '; Mode = 'AsIs' },
    @{ Cell = 'D28'; Value = '1.504'; Mode = 'AsText' },
    @{ Cell = 'E28'; Value = '  -1.83%  '; Mode = 'AsIs' },
    @{ Cell = 'E29'; Value = '  -0.38%  '; Mode = 'AsIs' },
    @{ Cell = 'D30'; Value = '4.120'; Mode = 'AsText' },
    @{ Cell = 'E30'; Value = '  -1.42%  '; Mode = 'AsIs' },
    @{ Cell = 'D31'; Value = '1.198'; Mode = 'AsText' },
    @{ Cell = 'E31'; Value = '  -0.22%  '; Mode = 'AsIs' },
    @{ Cell = 'D32'; Value = '0.05381'; Mode = 'AsText' },
    @{ Cell = 'E32'; Value = '  +5.52%  '; Mode = 'AsIs' },
    @{ Cell = 'E33'; Value = '  -2.02%  '; Mode = 'AsIs' },
    @{ Cell = 'D34'; Value = '0.7491'; Mode = 'AsText' },
    @{ Cell = 'E35'; Value = '  -2.44%  '; Mode = 'AsIs' },
    @{ Cell = 'E36'; Value = '  -0.06%  '; Mode = 'AsIs' },
    @{ Cell = 'D37'; Value = '1.298.54'; Mode = 'AsIs' },
    @{ Cell = 'E37'; Value = '  -2.76%  '; Mode = 'AsIs' },
    @{ Cell = 'D38'; Value = '0.01805'; Mode = 'AsText' },
    @{ Cell = 'E38'; Value = '  -3.54%  '; Mode = 'AsIs' },
    @{ Cell = 'D39'; Value = '2.707'; Mode = 'AsText' },
    @{ Cell = 'E39'; Value = '  -0.68%  '; Mode = 'AsIs' },
    @{ Cell = 'D40'; Value = '0.9307'; Mode = 'AsText' },
    @{ Cell = 'E40'; Value = '  -4.23%  '; Mode = 'AsIs' },
    @{ Cell = 'D41'; Value = '6.086'; Mode = 'AsText' },
    @{ Cell = 'E41'; Value = '  +4.74%  '; Mode = 'AsIs' },
    @{ Cell = 'D42'; Value = '0.00000000135'; Mode = 'AsText' },
    @{ Cell = 'E42'; Value = '  +9.38%  '; Mode = 'AsIs' },
    @{ Cell = 'D43'; Value = '0.08344'; Mode = 'AsText' },
    @{ Cell = 'E43'; Value = '  +28.67%  '; Mode = 'AsIs' },
    @{ Cell = 'D44'; Value = '104.38'; Mode = 'AsText' },
    @{ Cell = 'E44'; Value = '  -2.01%  '; Mode = 'AsIs' },
    @{ Cell = 'D45'; Value = '0.9989'; Mode = 'AsText' },
    @{ Cell = 'E45'; Value = '  -0.06%  '; Mode = 'AsIs' },
    @{ Cell = 'D46'; Value = '1.972.54'; Mode = 'AsIs' },
    @{ Cell = 'E46'; Value = '  -1.41%  '; Mode = 'AsIs' },
    @{ Cell = 'D47'; Value = '0.5173'; Mode = 'AsText' },
    @{ Cell = 'E47'; Value = '  -0.82%  '; Mode = 'AsIs' },
    @{ Cell = 'D48'; Value = '9.435'; Mode = 'AsText' },
    @{ Cell = 'E48'; Value = '  -3.49%  '; Mode = 'AsIs' },
    @{ Cell = 'B49'; Value = 'Aave'; Mode = 'AsIs' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; Mode = 'AsIs' },
    @{ Cell = 'D49'; Value = '63.57'; Mode = 'AsText' },
    @{ Cell = 'E49'; Value = '  -0.13%  '; Mode = 'AsIs' },
    @{ Cell = 'B50'; Value = 'RenderToken'; Mode = 'AsIs' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Mode = 'AsIs' },
    @{ Cell = 'D50'; Value = '1.756'; Mode = 'AsText' },
    @{ Cell = 'E50'; Value = '  -1.31%  '; Mode = 'AsIs' },
    @{ Cell = 'D51'; Value = '0.05924'; Mode = 'AsText' },
    @{ Cell = 'E51'; Value = '  -0.08%  '; Mode = 'AsIs' }
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    if ($edit.Mode -eq 'AsText') {
        $cell.NumberFormat = "@"
        $cell.Value = $edit.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $edit.Value
    }
}
